$d = $word.ActiveDocument

# 1) Ativação date update (Créditos paragraph)
$d.Content.Find.Execute("Ativação: 01/01/2021", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2024", 2)

# 2) "Programa" section (Portuguese) - the SECOND occurrence of this sentence
#    (the first occurrence, under "Programa resumido", must stay untouched).
#    Target paragraph 14 specifically via the Paragraphs collection.
$pPrograma = $d.Paragraphs(14).Range
$pPrograma.Find.Execute(
    "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1.Conceito de Projeto 2. Abordagem Preditiva, Adaptativa (Incremental e Interativo) e Híbrida 3. Abordagens Soft e hard de Gestão de Projetos 4. Gestão do Escopo e Visão do Projeto 5. Gestão do Tempo do Projeto 6. Gestão do Custo do Projeto 7. Gestão da Qualidade do Projeto 8. Gestão da Recursos do Projeto 9. Gestão das Comunicações do Projeto 10. Gestão dos Riscos do Projeto 11. Gestão das Aquisições do Projeto 12. Gestão das Partes Interessadas do Projeto 13. Maturidade em Gerenciamento de Projetos 14. Gestão de Portfólio de Projeto 15. Competências em Gestão de Projetos 16. Gestão Ágil e Lean de Projetos.17. Desenvolvimento de uma proposta de projeto para criação de um novo produto/processo/empreendimento.18. Visita (viagem didática complementar) a uma empresa para conhecer e entender os aspectos relacionados à Gestão de Projetos.",
    2)

# 3) "Programa" section (English/italic) - the SECOND occurrence of this sentence
#    (the first occurrence, under "Programa resumido", must stay untouched).
$pProgramaEn = $d.Paragraphs(15).Range
$pProgramaEn.Find.Execute(
    "General considerations about project management, Project initiation, Planning and Management Plan, Monitoring and EvaluationStructure, Execution and Control.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Project Concept 2. Predictive, Adaptive (Incremental and Interactive) and Hybrid Approach 3. Soft and hard Project Management Approaches 4. Project Scope and Vision Management 5. Project Time Management 6. Project Cost Management Project 7. Project Quality Management 8. Project Resource Management 9. Project Communications Management 10. Project Risk Management 11. Project Procurement Management 12. Project Stakeholder Management 13. Project Management Maturity Projects 14. Project Portfolio Management 15. Project Management Skills 16. Agile and Lean Project Management. 17. Development of a project proposal to create a new product/process/enterprise.18. Visit (complementary educational trip) to a company to learn about and understand aspects related to Project Management.",
    2)

# 4) Método (teaching method) text
$d.Content.Find.Execute("Aulas expositivas. Trabalhos em grupo. Seminários. Palestras. Exercícios em sala de aula.", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.", 2)

# 5) Critério (assessment criteria) text
$d.Content.Find.Execute("Duas Provas com peso de 30% cada uma. Trabalhos em sala de aula com peso de 20% e Trabalho final com peso de 20%", $true, $false, $false, $false, $false, $true, 1, $false, "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas", 2)

# 6) Norma de recuperação (recovery rule) text
$d.Content.Find.Execute("Prova única", $true, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação", 2)
